$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily update: the previous "last row" (row 23) loses its special
# date-only number format and reverts to the standard datetime format
# used by all the other data rows.
$ws.Range("A23").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append today's data as the new last row (row 24).
$ws.Range("A24").Value = 45973
$ws.Range("B24").Value = 53
$ws.Range("C24").Value = 61
$ws.Range("D24").Value = 60

# The new last row's date cell gets the date-only number format,
# matching the pattern previously applied to row 23.
$ws.Range("A24").NumberFormat = "YYYY-MM-DD"
